$d = $word.ActiveDocument
$brk = [char]11   # manual line break (w:br) marker for Range.Text assignment

function Replace-WithBreaks($findText, $newText) {
    $rng = $d.Content
    $rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    if (-not $rng.Find.Found) {
        Write-Output "NOT FOUND: $findText"
        return
    }
    $rng.Text = $newText
}

# --- "Programa resumido" section (Portuguese) ---
$old1 = "1.Cominuição e classificação de sólidos 2.Filtração3.Agitação de líquidos4.Trocadores de calor5.Evaporação6.Destilação7.Absorção8.Extração líquido-líquido"
$new1 = "1.Cominuição e classificação de sólidos $brk" + "2.Filtração$brk" + "3.Agitação de líquidos$brk" + "4.Trocadores de calor$brk" + "5.Evaporação$brk" + "6.Destilação$brk" + "7.Absorção$brk" + "8.Extração líquido-líquido"
Replace-WithBreaks $old1 $new1

# --- "Programa" section (Portuguese, identical text, second occurrence) ---
Replace-WithBreaks $old1 $new1

# --- Italic English program section ---
$old2 = "1. Comminution and solids classification2. Filtration3. Fluid mixing4. Heat exchangers5. Evaporation6. Distillation7. Absorption8. Liquid-liquid extraction"
$new2 = "1. Comminution and solids classification$brk" + "2. Filtration$brk" + "3. Fluid mixing$brk" + "4. Heat exchangers$brk" + "5. Evaporation$brk" + "6. Distillation$brk" + "7. Absorption$brk" + "8. Liquid-liquid extraction"
Replace-WithBreaks $old2 $new2

# --- Bibliografia section ---
$old3 = "1)COULSON, J. M.; RICHARDSON; J.F. Chemical Engineering. v. 2: Particle Technology e Separation Processes. 5ed. Amsterdan: Butterworth Heinemann, 1229p. 2005;2) COULSON & Richardson's Chemical Engineering: chemical engineering design by R.K. Sinnott. 6ed. Amsterdam: Elsevier Butterworth Heinemann, 895p. 2004;3)COUPER, J. R.; PENNEY, W. R.; FAIR, J. R.; W.; Stanley. M. Chemical Process Equipment: Selection and Design. 2ed. Amsterdam: Elsevier, 814p. 2005;4) FOUST, A. S.; WENZEL, L. A.; CLUMP, C. W.; MAUS, L.; ANDERSEN, L. B. 2ed. Princípios das operações unitárias. Rio de Janeiro: Guanabara Dois/LTC, 670p. 2008;5) GEANKOPLIS, C. J. Transport Processes and Separation Process Principles. 4ed. New York: Prentice Hall, 1026p. 2010;6) MCCABE, W. L.; SMITH, J. C.; HARRIOT, P. Unit operations of chemical engineering. 7ed. Boston: McGraw-Hill, 1140 p. 2005;7) PERRY's chemical engineers handbook. Editor in Chief Don W. Green; Late Editor Robert H. Perry New York: McGraw-Hill, 2008.8) SEADER, J. D; HENLEY, E. J. Separation Process Principles. 2ed. Hoboken, N.J: Wiley, 756p. 2006.9) TREYBAL, R. E. Mass-Transfer Operations. 3ed. Auckland: McGraw-Hill, 784p. 1980."
$new3 = "1)COULSON, J. M.; RICHARDSON; J.F. Chemical Engineering. v. 2: Particle Technology e Separation Processes. 5ed. Amsterdan: Butterworth Heinemann, 1229p. 2005;$brk" + "2) COULSON & Richardson's Chemical Engineering: chemical engineering design by R.K. Sinnott. 6ed. Amsterdam: Elsevier Butterworth Heinemann, 895p. 2004;$brk" + "3)COUPER, J. R.; PENNEY, W. R.; FAIR, J. R.; W.; Stanley. M. Chemical Process Equipment: Selection and Design. 2ed. Amsterdam: Elsevier, 814p. 2005;$brk" + "4) FOUST, A. S.; WENZEL, L. A.; CLUMP, C. W.; MAUS, L.; ANDERSEN, L. B. 2ed. Princípios das operações unitárias. Rio de Janeiro: Guanabara Dois/LTC, 670p. 2008;$brk" + "5) GEANKOPLIS, C. J. Transport Processes and Separation Process Principles. 4ed. New York: Prentice Hall, 1026p. 2010;$brk" + "6) MCCABE, W. L.; SMITH, J. C.; HARRIOT, P. Unit operations of chemical engineering. 7ed. Boston: McGraw-Hill, 1140 p. 2005;$brk" + "7) PERRY's chemical engineers handbook. Editor in Chief Don W. Green; Late Editor Robert H. Perry New York: McGraw-Hill, 2008.$brk" + "8) SEADER, J. D; HENLEY, E. J. Separation Process Principles. 2ed. Hoboken, N.J: Wiley, 756p. 2006.$brk" + "9) TREYBAL, R. E. Mass-Transfer Operations. 3ed. Auckland: McGraw-Hill, 784p. 1980."
Replace-WithBreaks $old3 $new3

Write-Output "Done"
